# Adds a "Gastos Mensal" / "Gastos annual" (Monthly/Annual Expenses) block
# of columns (H:K) to each of the four sections on Folha1, plus a totals
# row, mirroring the existing "Total mensal" / "Total anual" header style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header + data rows for each of the 4 sections, and the final totals row.
$headerRows = @(5, 9, 13, 17)
$dataRows   = @{ 6 = 400; 10 = 1100; 14 = 400; 18 = 500 }
$annualVals = @{ 6 = 48000; 10 = 13000; 14 = 4800; 18 = 6000 }

foreach ($r in $headerRows) {
    $hi = $ws.Range("H" + $r + ":I" + $r)
    $hi.Merge()
    $hi.Value = "Gastos Mensal"
    $hi.HorizontalAlignment = -4108
    $hi.Font.Bold = $true
    $hi.Font.Color = 16777215
    $hi.Interior.ColorIndex = 45
    $hi.Interior.Color = 12611584

    $jk = $ws.Range("J" + $r + ":K" + $r)
    $jk.Merge()
    $jk.Value = "Gastos annual"
    $jk.HorizontalAlignment = -4108
    $jk.Font.Bold = $true
    $jk.Font.Color = 16777215
    $jk.Interior.Color = 12611584
}

foreach ($r in $dataRows.Keys) {
    $hi = $ws.Range("H" + $r + ":I" + $r)
    $hi.Merge()
    $hi.Value = $dataRows[$r]
    $hi.HorizontalAlignment = -4108
    $hi.Font.Bold = $true
    $hi.Interior.Color = 10284031

    $jk = $ws.Range("J" + $r + ":K" + $r)
    $jk.Merge()
    $jk.Value = $annualVals[$r]
    $jk.HorizontalAlignment = -4108
    $jk.Interior.Color = 10284031
}

$hi20 = $ws.Range("H20:I20")
$hi20.Merge()
$hi20.Value = 6400
$hi20.HorizontalAlignment = -4108

$jk20 = $ws.Range("J20:K20")
$jk20.Merge()
$jk20.Value = 76680
$jk20.HorizontalAlignment = -4108

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "done"
